# Update records for two candidates on the "Registros" sheet:
# their score and ethnicity answers changed, plus answers for
# questions 4, 6, 8 and 19 were revised.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Registros")

# Row 7 - Maria Aparecida Silva
$ws.Range("B7").Value = "26"
$ws.Range("C7").Value = "no"
$ws.Range("G7").Value = "0"
$ws.Range("I7").Value = "0"
$ws.Range("K7").Value = "1"
$ws.Range("V7").Value = "1"

# Row 9 - João Pedro Santos
$ws.Range("B9").Value = "26"
$ws.Range("C9").Value = "no"
$ws.Range("G9").Value = "0"
$ws.Range("I9").Value = "0"
$ws.Range("K9").Value = "1"
$ws.Range("V9").Value = "1"
